# LogisticExpenseProductTemplate update
# - zoneArea changes from 1263 to 2910 (text, quote-prefixed like the original)
# - zoneClassPrice collapses to BKK1 for every row (was BKK1..BKK7)
# - effectiveDate moves from 2021-06-17 (44364) to 2023-09-08 (45177)
# - logisticExp increments 2.5 .. 5.5 in steps of 0.5 instead of a flat 2.5
# - productCode / productName / status columns are left untouched
# - column D is widened and the active selection moves to D7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newZoneArea  = "'2910"
$newZoneClass = "BKK1"
$newDate      = 45177

$logisticExp = @(2.5, 3, 3.5, 4, 4.5, 5, 5.5)

for ($i = 0; $i -lt 7; $i++) {
    $row = 2 + $i

    $ws.Cells.Item($row, 1).Value = $newZoneArea
    $ws.Cells.Item($row, 2).Value = $newZoneClass
    $ws.Cells.Item($row, 3).Value = $newDate
    $ws.Cells.Item($row, 6).Value = $logisticExp[$i]
}

# Widen column D and move the selection to D7, matching the refreshed layout
$ws.Columns.Item(4).ColumnWidth = 20.76
$ws.Range("D7").Select() | Out-Null
